$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Driver Summary")

# Bad Drivers table (row 3 is the single driver row, row 4 is Totals)
$ws.Range("C3").Value = 121
$ws.Range("D3").Value = 95.8
$ws.Range("C4").Value = 121

# Good Drivers table - update total samples for the "23.20.1.1" driver row (row 14)
$ws.Range("B14").Value = 15827
